$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.954.21"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.894.01"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.8342"
$ws.Range("E5").Value = "  +8.55%  "
$ws.Range("D6").Value = "'241.84"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.3242"
$ws.Range("E8").Value = "  +6.35%  "
$ws.Range("D9").Value = "'26.74"
$ws.Range("E9").Value = "  +5.57%  "
$ws.Range("D10").Value = "'0.07032"
$ws.Range("D11").Value = "'0.08041"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "'0.7487"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "1.898.94"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "'5.208"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'92.42"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "29.954.35"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "'14.08"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").Value = "'5.921"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'244.56"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "'0.000007756"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D22").Value = "2.148.61"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'6.921"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "'0.1616"
$ws.Range("E25").Value = "  +25.57%  "
$ws.Range("D26").Value = "'168.04"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").Value = "'9.200"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'18.87"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "'4.268"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'0.05676"
$ws.Range("E33").Value = "  +7.49%  "
$ws.Range("D34").Value = "'4.080"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "'1.280"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("D36").Value = "'0.7344"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "'2.719"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.01911"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'2.779"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'0.4424"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "'72.07"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "'5.954"
$ws.Range("E42").Value = "  -3.78%  "
$ws.Range("D43").Value = "'0.8434"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D45").Value = "'1.892"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "'101.05"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").Value = "'7.597"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.696"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'991.49"
$ws.Range("E49").Value = "  +9.22%  "
$ws.Range("D50").Value = "2.049.81"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "'36.06"
$ws.Range("E51").Value = "  -0.17%  "
